$wb = $excel.ActiveWorkbook

# --- mock1 sheet updates ---
$ws1 = $wb.Worksheets.Item("mock1")

$ws1.Range("C2").Value = 8864
$ws1.Range("D2").Value = 8865
$ws1.Range("G2").Value = 1726

$ws1.Range("C3").Value = 2283
$ws1.Range("D3").Value = 2283

$ws1.Range("C5").Value = 1233
$ws1.Range("D5").Value = 1234

$ws1.Range("C6").Value = 931
$ws1.Range("D6").Value = 931
$ws1.Range("G6").Value = 258

$ws1.Range("C7").Value = 808
$ws1.Range("D7").Value = 808

$ws1.Range("G9").Value = 199

$ws1.Range("C10").Value = 158
$ws1.Range("D10").Value = 158

# --- mock2 sheet updates ---
$ws2 = $wb.Worksheets.Item("mock2")

$ws2.Range("C2").Value = 4638
$ws2.Range("D2").Value = 4630
$ws2.Range("E2").Value = 877

$ws2.Range("C3").Value = 1961
$ws2.Range("D3").Value = 1961

$ws2.Range("C4").Value = 1494
$ws2.Range("D4").Value = 1494

$ws2.Range("C5").Value = 617
$ws2.Range("D5").Value = 617

$ws2.Range("E6").Value = 44
